$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "SCD0022"

# Update TC_ID column (B) values from DGS-329 to SCD0022-001
$ws.Range("B2").Value = "SCD0022-001"
$ws.Range("B3").Value = "SCD0022-001"
$ws.Range("B4").Value = "SCD0022-001"
$ws.Range("B5").Value = "SCD0022-001"

# Set column B width to reflect new content width (~12.86 chars, bestFit for "SCD0022-001")
$ws.Columns.Item(2).ColumnWidth = 12

# Update selection to reflect final cursor position, and zoom level
$ws.Select()
$ws.Range("B6").Select()
$excel.ActiveWindow.Zoom = 66
